$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.334.78'
$ws.Range("E2").Value = '  +3.28%  '
$ws.Range("D3").Value = '1.723.15'
$ws.Range("E3").Value = '  +3.29%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9988'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9993'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4735'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2640'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06196'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("D10").Value = '1.718.84'
$ws.Range("E10").Value = '  +3.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07064'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.45'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5970'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.430'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.23'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9994'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '26.334.14'
$ws.Range("E17").Value = '  +3.29%  '
$ws.Range("B18").Value = 'BinanceUSD'
$ws.Range("C18").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9998'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006827'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.04%  '
$ws.Range("E20").Value = '  +0.94%  '
$ws.Range("D21").Value = '1.938.82'
$ws.Range("E21").Value = '  +3.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.527'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.734'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").Value = '  -0.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '135.08'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.23%  '
$ws.Range("E26").Value = '  +1.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.768'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.400'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.961'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.685'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.12%  '
$ws.Range("E32").Value = '  +0.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04500'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.96%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.613'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9835'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.74%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6239'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.32%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9273'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.84%  '
$ws.Range("B38").Value = 'Quant'
$ws.Range("C38").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '113.58'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +16.92%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.453'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.48%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.933'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.80%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.666'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +17.25%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01488'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.08%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3834'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.99%  '
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1187'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.12%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.358'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.57%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05272'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.31%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.878'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.49%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.34%  '
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3388'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.91%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.219'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.10%  '

Write-Host "Update complete"